$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 271; existing rows 271:286 shift down to 272:287.
$ws.Rows("271:271").Insert()

# Fill the newly inserted row 271 with the new data record.
$ws.Cells.Item(271, 1).Value = 10
$ws.Cells.Item(271, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(271, 3).Value = "La Araucanía"
$ws.Cells.Item(271, 4).Value = 44706
$ws.Cells.Item(271, 4).NumberFormat = $ws.Cells.Item(272, 4).NumberFormat
$ws.Cells.Item(271, 5).Value = 9
$ws.Cells.Item(271, 6).Value = 100112001
$ws.Cells.Item(271, 7).Value = "Berenjena"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 50
$ws.Cells.Item(271, 11).Value = 10000
$ws.Cells.Item(271, 12).Value = 10000
$ws.Cells.Item(271, 13).Value = 10000
$ws.Cells.Item(271, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(271, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(271, 16).Value = 167
$ws.Cells.Item(271, 17).Value = 60
$ws.Cells.Item(271, 18).Value = "Hortaliza"
